# Update Branch - Brand - SKU wise Stock Aging Status: Summary
# Columns: D=Super Under Stock, E=Under Stock, F=Normal Stock, G=Over Stock, H=Super Over Stock

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(34, 14, 6, 13, 39)
    3  = @(37, 16, 7, 9, 37)
    4  = @(38, 15, 7, 11, 35)
    5  = @(38, 11, 5, 14, 38)
    6  = @(35, 6, 13, 15, 37)
    7  = @(31, 15, 4, 10, 46)
    8  = @(38, 9, 11, 6, 42)
    9  = @(32, 7, 5, 12, 50)
    10 = @(35, 16, 6, 14, 35)
    11 = @(37, 6, 7, 10, 46)
    12 = @(36, 9, 10, 10, 41)
    13 = @(37, 19, 8, 11, 31)
    14 = @(35, 15, 6, 6, 44)
    15 = @(30, 11, 4, 9, 52)
    16 = @(36, 12, 8, 11, 39)
    17 = @(41, 21, 4, 9, 31)
    18 = @(31, 7, 3, 8, 57)
    19 = @(34, 8, 2, 5, 57)
    20 = @(35, 13, 13, 14, 31)
    21 = @(35, 11, 4, 4, 52)
    22 = @(39, 12, 13, 8, 34)
    23 = @(42, 17, 8, 11, 28)
    24 = @(45, 11, 4, 9, 37)
    25 = @(31, 3, 4, 6, 62)
    26 = @(34, 13, 7, 9, 43)
    27 = @(31, 14, 6, 8, 47)
    28 = @(37, 11, 5, 12, 41)
    29 = @(35, 7, 4, 10, 50)
    30 = @(39, 15, 11, 13, 28)
    31 = @(38, 14, 8, 9, 37)
    32 = @(38, 10, 7, 7, 44)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
    $ws.Cells.Item($row, 6).Value = $vals[2]
    $ws.Cells.Item($row, 7).Value = $vals[3]
    $ws.Cells.Item($row, 8).Value = $vals[4]
}
